$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 3 (NganLNT -> QuangTrung) ---
$ws.Range("B3").Value = "QuangTrung"
$ws.Range("C3").Value = "Trần Quang Trung"
$ws.Range("E3").Value = "quangtrung@gmail.com"
$ws.Range("F3").Value = "data/face_train/QuangTrung/2021-05-17-14-12-56-012379.jpg"

# --- Update row 2's avatar path with the new placeholder text ---
$ws.Range("F2").Value = "sdafdsafasdfsad"

# --- Remove the now unused placeholder rows 4-10 ---
$ws.Rows("4:10").Delete()

# --- Match the selection left by the author after the edit ---
[void]$ws.Range("F2").Select()
